$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 2 - this shifts the existing
# rows 2..9 down to 3..10 (and bumps the sheet dimension to A1:N10).
$ws.Rows.Item(2).Insert()

# Fill in the new row 2 with the new archive entry.
# Date-looking text (yyyy-mm-dd) is prefixed with a leading apostrophe so
# Excel's type inference stores it as literal text instead of auto-
# converting it to a date serial number. A lone apostrophe likewise forces
# an (empty) text cell for L2 instead of clearing the cell outright.
$ws.Range("A2").Value = "'2025-08-06"
$ws.Range("B2").Value = "12:16:42"
$ws.Range("C2").Value = "cibul"
$ws.Range("D2").Value = "Položka"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "sedák"
$ws.Range("G2").Value = "EU-SVA-999999-25"
$ws.Range("H2").Value = "BL5"
$ws.Range("I2").Value = "2 ks"
$ws.Range("J2").Value = "'2025-08-01"
$ws.Range("K2").Value = "Expirace"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'2025-08-06"
$ws.Range("N2").Value = "Cibulka"

# The inserted row inherited the bold/bordered header formatting (Excel's
# "format from row above" behaviour on row insert) and the apostrophe
# prefixes set a quote-prefix style flag; reset the whole row back to the
# plain "Normal" style used by every other data row.
$ws.Range("A2:N2").Style = "Normal"
